$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-98 from 45188 to 45189
for ($r = 2; $r -le 98; $r++) {
    $ws.Cells.Item($r, 3).Value = 45189
}
